# Applies the "Updated symbol list" edit: refreshed prices/volumes and a
# re-ranking of several coins (rows 10-18, 41-43) in the cryptos sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cell -> new value, in the same order as the source diff.
$updates = @(
    @{ Cell = "D2"; Value = "246.31" }
    @{ Cell = "D3"; Value = "23.97" }
    @{ Cell = "D4"; Value = "5.359" }
    @{ Cell = "D5"; Value = "0.05819" }
    @{ Cell = "D6"; Value = "3.370" }
    @{ Cell = "D7"; Value = "6.477" }
    @{ Cell = "D8"; Value = "0.8102" }
    @{ Cell = "D9"; Value = "0.9217" }
    @{ Cell = "B10"; Value = "One" }
    @{ Cell = "C10"; Value = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one" }
    @{ Cell = "D10"; Value = "0.0005965" }
    @{ Cell = "E10"; Value = "9OneONE" }
    @{ Cell = "B11"; Value = "WazirX" }
    @{ Cell = "C11"; Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx" }
    @{ Cell = "D11"; Value = "0.1401" }
    @{ Cell = "E11"; Value = "10WazirXWRX" }
    @{ Cell = "B12"; Value = "MandalaExchangeToken" }
    @{ Cell = "C12"; Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx" }
    @{ Cell = "D12"; Value = "0.07397" }
    @{ Cell = "E12"; Value = "11MandalaExchangeTokenMDX" }
    @{ Cell = "B13"; Value = "LiechtensteinCryptoassetsExchange" }
    @{ Cell = "C13"; Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx" }
    @{ Cell = "D13"; Value = "0.03208" }
    @{ Cell = "E13"; Value = "12LiechtensteinCryptoassetsExchangeLCX" }
    @{ Cell = "B14"; Value = "BitrueCoin" }
    @{ Cell = "C14"; Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr" }
    @{ Cell = "D14"; Value = "0.03021" }
    @{ Cell = "E14"; Value = "13BitrueCoinBTR" }
    @{ Cell = "B15"; Value = "BitMartToken" }
    @{ Cell = "C15"; Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx" }
    @{ Cell = "D15"; Value = "0.09384" }
    @{ Cell = "E15"; Value = "14BitMartTokenBMX" }
    @{ Cell = "B16"; Value = "MCDex" }
    @{ Cell = "C16"; Value = "https://coinranking.com/coin/3nMM61qeg+mcdex-mcb" }
    @{ Cell = "D16"; Value = "3.853" }
    @{ Cell = "E16"; Value = "15MCDexMCB" }
    @{ Cell = "B17"; Value = "BitForexToken" }
    @{ Cell = "C17"; Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf" }
    @{ Cell = "D17"; Value = "0.001555" }
    @{ Cell = "E17"; Value = "16BitForexTokenBF" }
    @{ Cell = "B18"; Value = "CoinExToken" }
    @{ Cell = "C18"; Value = "https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet" }
    @{ Cell = "D18"; Value = "0.04702" }
    @{ Cell = "E18"; Value = "17CoinExTokenCET" }
    @{ Cell = "D19"; Value = "0.006019" }
    @{ Cell = "D20"; Value = "0.001247" }
    @{ Cell = "D22"; Value = "0.00008793" }
    @{ Cell = "E22"; Value = "21NitroExNTX" }
    @{ Cell = "D23"; Value = "3.595" }
    @{ Cell = "D26"; Value = "0.1318" }
    @{ Cell = "D28"; Value = "0.0002348" }
    @{ Cell = "D41"; Value = "0.006405" }
    @{ Cell = "E41"; Value = "40KickTokenKICKBestin24h" }
    @{ Cell = "B42"; Value = "CEJI" }
    @{ Cell = "C42"; Value = "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji" }
    @{ Cell = "D42"; Value = "0.003497" }
    @{ Cell = "E42"; Value = "41CEJICEJI" }
    @{ Cell = "B43"; Value = "BKEXToken" }
    @{ Cell = "C43"; Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk" }
    @{ Cell = "D43"; Value = "0.1067" }
    @{ Cell = "E43"; Value = "42BKEXTokenBKK" }
    @{ Cell = "D44"; Value = "0.008625" }
    @{ Cell = "D45"; Value = "0.00005268" }
    @{ Cell = "D46"; Value = "0.00000000749" }
    @{ Cell = "D47"; Value = "0.7094" }
    @{ Cell = "D48"; Value = "0.001836" }
    @{ Cell = "D49"; Value = "0.00002098" }
    @{ Cell = "D50"; Value = "0.0001998" }
)

foreach ($u in $updates) {
    $range = $ws.Range($u.Cell)
    if ($u.Cell -match "^D\d+$") {
        # Column D stores prices as plain text in this sheet (t="inlineStr").
        # Force a text format before assigning so Excel does not silently
        # convert the numeric-looking string into a real number, then drop
        # the format override again so the cell keeps its original (default)
        # style, matching the source workbook exactly.
        $range.NumberFormat = "@"
        $range.Value = $u.Value
        $range.Style = "Normal"
    } else {
        $range.Value = $u.Value
    }
}
